$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$notes = @("C","C#","D","D#","E","F","F#","G","G#","A","A#","B")

function Get-MidiNoteName($i) {
    $octave = [Math]::Floor($i / 12) - 2
    $name = $notes[$i % 12]
    return "$name$octave ($i)"
}

# Sheet2 "DO NOT MODIFY!" column A, rows 3..130 hold the 128 MIDI note
# names (MIDI note 0..127) in the old "NNN:NOTE" format. Rewrite them to
# the new "NOTE (NNN)" format.
for ($i = 0; $i -lt 128; $i++) {
    $row = $i + 3
    $ws2.Cells.Item($row, 1).Value = (Get-MidiNoteName $i)
}

# Sheet1 "VHG Mono (Main)" column E, rows 3..22 hold a copy of the note
# name (MIDI notes 9..28) picked via the data-validation list sourced
# from sheet2 column A. Update them to match the new format too.
for ($i = 9; $i -le 28; $i++) {
    $row = ($i - 9) + 3
    $ws1.Cells.Item($row, 5).Value = (Get-MidiNoteName $i)
}

# Update the frozen-pane selection on sheet1: bottomRight pane active
# cell moves from E3 to I2.
$ws1.Range("I2").Select()
